$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Roster20171031")

# ---------------------------------------------------------------------
# 1. Fill in sample pilot data on row 2 (A2:M2).
#    Most of these cells already carry the text-style (s=1) from the
#    template, so a direct .Value assignment keeps that style.
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "Pilot"
$ws.Range("B2").Value = "Test"
$ws.Range("C2").Value = "User"
$ws.Range("D2").Value = "TUP"

# E2 has no pre-existing style in the template (the row skips E2), so
# give it the same "text" number format as its neighbours before
# writing the value - otherwise the new cell would default to style 0.
$ws.Range("E2").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("E2").Value = "123 Main St"

$ws.Range("F2").Value = "Somewhere, AL 12345"
$ws.Range("G2").Value = "555 123 4567"
$ws.Range("H2").Value = "123 555 7890"

# I2 gets the e-mail address plus a mailto: hyperlink (adds the
# "Hyperlink" cell style/font automatically).
$ws.Range("I2").Value = "example@domain.com"
$ws.Hyperlinks.Add($ws.Range("I2"), "mailto:example@domain.com") | Out-Null

$ws.Range("J2").Value = "PDK"

# K2 ("1/1/1901") must stay a literal text string, not turn into a
# date serial number. Writing it through a TEXT() formula elsewhere
# and pasting the computed value back (values-only) keeps it a plain
# string without disturbing the cell's existing style.
$ws.Range("Z1").Formula = '=TEXT("1/1/1901","M/D/YYYY")'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("K2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues) | Out-Null
$ws.Range("Z1").ClearContents()

$ws.Range("L2").Value = "USA"
$ws.Range("M2").Value = "SIC"

# ---------------------------------------------------------------------
# 2. Row 2's helper-column formulas stay, but the SamAccountName
#    formula in R2 now produces the shorter suffixes.
# ---------------------------------------------------------------------
$ws.Range("R2").Formula = '=IF(M2="PIC",IF(J2="PSM","pc12picpsm",IF(Q2="North","PC12PICRemoteN","PC12PICRemoteS")),IF(J2="PSM","PC12SICPSM",IF(Q2="North","PC12SICRemoteN","PC12SICRemoteS")))'

# ---------------------------------------------------------------------
# 3. Rows 3-11 no longer carry the (previously erroring) helper
#    formulas - the log-entry-for-groups move means the template rows
#    are just blank, styled cells now.
# ---------------------------------------------------------------------
$ws.Range("N3:R11").ClearContents()

# ---------------------------------------------------------------------
# 4. Selection / view bookkeeping - select the data block and drop the
#    scrolled-right "topLeftCell".
# ---------------------------------------------------------------------
$ws.Range("A1").Select() | Out-Null
$ws.Range("A3:R11").Select() | Out-Null
